$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 8 (hunk 0)
$ws.Range("H8").Value = 733.3333
$ws.Range("I8").Value = 100
$ws.Range("J8").Value = 2000
$ws.Range("K8").Value = 300
$ws.Range("L8").Value = 6000
$ws.Range("M8").Value = -161
$ws.Range("N8").Value = -6278

# row 11 (hunk 1)
$ws.Range("H11").Value = 630
$ws.Range("I11").Value = 630
$ws.Range("K11").Value = 630
$ws.Range("M11").Value = -490

# row 113 (hunk 2)
$ws.Range("H113").Value = 1901.2222
$ws.Range("I113").Value = 1684.1666
$ws.Range("J113").Value = 2335.3333
$ws.Range("K113").Value = 1684.1666
$ws.Range("L113").Value = 2335.3333
$ws.Range("M113").Value = 1569.8334
$ws.Range("N113").Value = -8843.3333

# row 138 (hunk 3)
$ws.Range("H138").Value = 3204.26
$ws.Range("I138").Value = 1716.5416
$ws.Range("J138").Value = 3674.0657
$ws.Range("K138").Value = 5149.6248
$ws.Range("L138").Value = 11022.1971
$ws.Range("M138").Value = -9.624799999999595
$ws.Range("N138").Value = -21302.1971

$ws = $wb.Worksheets.Item("ARM")
# row 74 (hunk 4)
$ws.Range("H74").Value = 25129.934
$ws.Range("I74").Value = 2661.8215
$ws.Range("K74").Value = 2661.8215
$ws.Range("M74").Value = -1787.8215

# row 77 (hunk 5)
$ws.Range("H77").Value = 25129.934
$ws.Range("I77").Value = 2661.8215
$ws.Range("K77").Value = 13309.1075
$ws.Range("M77").Value = -8941.1075

$ws = $wb.Worksheets.Item("BSM")
# row 12 (hunk 6)
$ws.Range("H12").Value = 4469.6
$ws.Range("I12").Value = 637
$ws.Range("J12").Value = 19800
$ws.Range("K12").Value = 637
$ws.Range("L12").Value = 19800
$ws.Range("M12").Value = -469
$ws.Range("N12").Value = -20136

$ws = $wb.Worksheets.Item("CRP")
# row 31 (hunk 7)
$ws.Range("H31").Value = 3563.1482
$ws.Range("I31").Value = 1773.3334
$ws.Range("J31").Value = 7142.778
$ws.Range("K31").Value = 1773.3334
$ws.Range("L31").Value = 7142.778
$ws.Range("M31").Value = -1478.3334
$ws.Range("N31").Value = -7732.778

# row 34 (hunk 8)
$ws.Range("H34").Value = 3563.1482
$ws.Range("I34").Value = 1773.3334
$ws.Range("J34").Value = 7142.778
$ws.Range("K34").Value = 1773.3334
$ws.Range("L34").Value = 7142.778
$ws.Range("M34").Value = -1571.3334
$ws.Range("N34").Value = -7546.778

# row 141 (hunk 9)
$ws.Range("H141").Value = 45143.89
$ws.Range("I141").Value = 20296
$ws.Range("J141").Value = 48249.875
$ws.Range("K141").Value = 20296
$ws.Range("L141").Value = 48249.875
$ws.Range("M141").Value = -15116
$ws.Range("N141").Value = -58609.875

$ws = $wb.Worksheets.Item("CUL")
# row 4 (hunk 10)
$ws.Range("H4").Value = 100.97727
$ws.Range("I4").Value = 99.02631
$ws.Range("K4").Value = 297.07893
$ws.Range("M4").Value = -185.07893

# row 58 (hunk 11)
$ws.Range("H58").Value = 3125.8572
$ws.Range("I58").Value = 1166.6666
$ws.Range("J58").Value = 3452.389
$ws.Range("K58").Value = 3499.9998
$ws.Range("L58").Value = 10357.167
$ws.Range("M58").Value = -3371.9998
$ws.Range("N58").Value = -10613.167

# row 62 (hunk 12)
$ws.Range("H62").Value = 570
$ws.Range("I62").Value = 570
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 1710
$ws.Range("L62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -1024

# row 63 (hunk 13)
$ws.Range("H63").Value = 4775.294
$ws.Range("I63").Value = 2297.7778
$ws.Range("J63").Value = 7562.5
$ws.Range("K63").Value = 6893.3334
$ws.Range("L63").Value = 22687.5
$ws.Range("M63").Value = -6144.3334
$ws.Range("N63").Value = -24185.5

# row 65 (hunk 14)
$ws.Range("H65").Value = 570
$ws.Range("I65").Value = 570
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 5130
$ws.Range("L65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -1698

# row 66 (hunk 15)
$ws.Range("H66").Value = 4775.294
$ws.Range("I66").Value = 2297.7778
$ws.Range("J66").Value = 7562.5
$ws.Range("K66").Value = 20680.0002
$ws.Range("L66").Value = 68062.5
$ws.Range("M66").Value = -16936.0002
$ws.Range("N66").Value = -75550.5

# row 68 (hunk 16)
$ws.Range("H68").Value = 935.8
$ws.Range("I68").Value = 590
$ws.Range("J68").Value = 1166.3334
$ws.Range("K68").Value = 1770
$ws.Range("L68").Value = 3499.0002
$ws.Range("M68").Value = -959
$ws.Range("N68").Value = -5121.0002

# row 71 (hunk 17)
$ws.Range("H71").Value = 935.8
$ws.Range("I71").Value = 590
$ws.Range("J71").Value = 1166.3334
$ws.Range("K71").Value = 5310
$ws.Range("L71").Value = 10497.0006
$ws.Range("M71").Value = -1254
$ws.Range("N71").Value = -18609.0006

# row 75 (hunk 18)
$ws.Range("H75").Value = 5337
$ws.Range("I75").Value = 2222
$ws.Range("J75").Value = 5960
$ws.Range("K75").Value = 6666
$ws.Range("L75").Value = 17880
$ws.Range("M75").Value = -5668
$ws.Range("N75").Value = -19876

# row 78 (hunk 19)
$ws.Range("H78").Value = 5337
$ws.Range("I78").Value = 2222
$ws.Range("J78").Value = 5960
$ws.Range("K78").Value = 19998
$ws.Range("L78").Value = 53640
$ws.Range("M78").Value = -15006
$ws.Range("N78").Value = -63624

# row 86 (hunk 20)
$ws.Range("H86").Value = 1497.5
$ws.Range("J86").Value = 1497.5
$ws.Range("L86").Value = 4492.5
$ws.Range("N86").Value = -6864.5

# row 87 (hunk 21)
$ws.Range("H87").Value = 4750
$ws.Range("I87").Value = 4750
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 14250
$ws.Range("L87").Value = 0
$ws.Range("M87").ClearContents()
$ws.Range("N87").Value = -13002

# row 89 (hunk 22)
$ws.Range("H89").Value = 1497.5
$ws.Range("J89").Value = 1497.5
$ws.Range("L89").Value = 13477.5
$ws.Range("N89").Value = -25333.5

# row 90 (hunk 23)
$ws.Range("H90").Value = 4750
$ws.Range("I90").Value = 4750
$ws.Range("J90").Value = 0
$ws.Range("K90").Value = 42750
$ws.Range("L90").Value = 0
$ws.Range("M90").ClearContents()
$ws.Range("N90").Value = -36510

# row 120 (hunk 24)
$ws.Range("H120").Value = 10141.111
$ws.Range("I120").Value = 3748.889
$ws.Range("J120").Value = 16533.334
$ws.Range("K120").Value = 11246.667
$ws.Range("L120").Value = 49600.00199999999
$ws.Range("M120").Value = -6408.667000000001
$ws.Range("N120").Value = -59276.00199999999

$ws = $wb.Worksheets.Item("GSM")
# row 41 (hunk 25)
$ws.Range("H41").Value = 3650
$ws.Range("I41").Value = 475
$ws.Range("J41").Value = 10000
$ws.Range("K41").Value = 475
$ws.Range("L41").Value = 10000
$ws.Range("M41").Value = -120
$ws.Range("N41").Value = -10710

# row 113 (hunk 26)
$ws.Range("H113").Value = 2026.1538
$ws.Range("I113").Value = 1886
$ws.Range("K113").Value = 1886
$ws.Range("M113").Value = 284

$ws = $wb.Worksheets.Item("LTW")
# row 7 (hunk 27)
$ws.Range("H7").Value = 7836870.5
$ws.Range("I7").Value = 12784302
$ws.Range("J7").Value = 3436.1667
$ws.Range("K7").Value = 12784302
$ws.Range("L7").Value = 3436.1667
$ws.Range("M7").Value = -12784190
$ws.Range("N7").Value = -3660.1667

# row 22 (hunk 28)
$ws.Range("H22").Value = 1174.8334
$ws.Range("I22").Value = 979.6
$ws.Range("J22").Value = 1314.2858
$ws.Range("K22").Value = 979.6
$ws.Range("L22").Value = 1314.2858
$ws.Range("M22").Value = -684.6
$ws.Range("N22").Value = -1904.2858

# row 27 (hunk 29)
$ws.Range("H27").Value = 1174.8334
$ws.Range("I27").Value = 979.6
$ws.Range("J27").Value = 1314.2858
$ws.Range("K27").Value = 979.6
$ws.Range("L27").Value = 1314.2858
$ws.Range("M27").Value = -872.6
$ws.Range("N27").Value = -1528.2858

# row 46 (hunk 30)
$ws.Range("H46").Value = 2240
$ws.Range("I46").Value = 1316.6666
$ws.Range("J46").Value = 3625
$ws.Range("K46").Value = 1316.6666
$ws.Range("L46").Value = 3625
$ws.Range("M46").Value = -1128.6666
$ws.Range("N46").Value = -4001

# row 61 (hunk 31)
$ws.Range("H61").Value = 4285.7144
$ws.Range("I61").Value = 2666.6667
$ws.Range("J61").Value = 5500
$ws.Range("K61").Value = 2666.6667
$ws.Range("L61").Value = 5500
$ws.Range("M61").Value = -2464.6667
$ws.Range("N61").Value = -5904

# row 113 (hunk 32)
$ws.Range("H113").Value = 4285.7144
$ws.Range("I113").Value = 2666.6667
$ws.Range("J113").Value = 5500
$ws.Range("K113").Value = 2666.6667
$ws.Range("L113").Value = 5500
$ws.Range("M113").Value = -496.6667000000002
$ws.Range("N113").Value = -9840

# row 126 (hunk 33)
$ws.Range("H126").Value = 7836870.5
$ws.Range("I126").Value = 12784302
$ws.Range("J126").Value = 3436.1667
$ws.Range("K126").Value = 38352906
$ws.Range("L126").Value = 10308.5001
$ws.Range("M126").Value = -38350436
$ws.Range("N126").Value = -15248.5001

# row 132 (hunk 34)
$ws.Range("H132").Value = 2938.875
$ws.Range("J132").Value = 4423.75
$ws.Range("L132").Value = 13271.25
$ws.Range("N132").Value = -18331.25
